$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.337.30"
$ws.Range("E2").Value = "  -4.57%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.313.94"
$ws.Range("E3").Value = "  -5.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.44"
$ws.Range("E5").Value = "  -3.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.93"
$ws.Range("E6").Value = "  -2.88%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.313.90"
$ws.Range("E8").Value = "  -5.14%  "

# Row 9
$ws.Range("E9").Value = "  -0.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("E10").Value = "  -5.03%  "

# Row 11
$ws.Range("E11").Value = "  -4.07%  "

# Row 12
$ws.Range("E12").Value = "  -2.77%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.881.93"
$ws.Range("E13").Value = "  -5.07%  "

# Row 14
$ws.Range("E14").Value = "  -0.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.313.52"
$ws.Range("E15").Value = "  -5.11%  "

# Row 16
$ws.Range("E16").Value = "  -5.34%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.60"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.422.19"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.66"
$ws.Range("E19").Value = "  -1.55%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("E20").Value = "  -0.92%  "

# Row 21
$ws.Range("E21").Value = "  -10.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "355.72"
$ws.Range("E22").Value = "  -7.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  -3.70%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.447.91"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.43"
$ws.Range("E26").Value = "  -6.67%  "

# Row 27
$ws.Range("E27").Value = "  -5.43%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("E30").Value = "  -1.23%  "

# Row 31
$ws.Range("E31").Value = "  -1.51%  "

# Row 32
$ws.Range("E32").Value = "  -5.81%  "

# Row 33
$ws.Range("E33").Value = "  -2.91%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.343.81"
$ws.Range("E35").Value = "  -5.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.60"
$ws.Range("E36").Value = "  -2.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("E37").Value = "  -0.57%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.81"
$ws.Range("E38").Value = "  -0.60%  "

# Row 39: 'ImmutableX' -> 'Monero'
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "161.14"
$ws.Range("E39").Value = "  -1.77%  "

# Row 40: 'Monero' -> 'ImmutableX'
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  -2.86%  "

# Row 41
$ws.Range("E41").Value = "  -2.78%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
$ws.Range("E43").Value = "  +0.59%  "

# Row 45
$ws.Range("E45").Value = "  -7.56%  "

# Row 46
$ws.Range("E46").Value = "  -4.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.56"
$ws.Range("E47").Value = "  -4.78%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.41"
$ws.Range("E48").Value = "  -7.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.71"

# Row 50: 'SuiNetwork' -> 'InjectiveProtocol'
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.41"
$ws.Range("E50").Value = "  +3.32%  "

# Row 51: 'InjectiveProtocol' -> 'SuiNetwork'
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.857"
$ws.Range("E51").Value = "  -8.56%  "

